$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.173414826393127
$ws.Range("B1").Value = 2.436054944992065
$ws.Range("D1").Value = 2.366298913955688
$ws.Range("E1").Value = 1.237086176872253
